$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns F,G,H,I left by one position into F,G,H for every row
# (header row 1 through last data row 246), then drop the now-redundant
# column I. Column D's header additionally gets a new "ppm" label.
for ($r = 1; $r -le 246; $r++) {
    $g = $ws.Cells.Item($r, 7).Value()
    $h = $ws.Cells.Item($r, 8).Value()
    $i = $ws.Cells.Item($r, 9).Value()
    $ws.Cells.Item($r, 6).Value = $g
    $ws.Cells.Item($r, 7).Value = $h
    $ws.Cells.Item($r, 8).Value = $i
}

$ws.Range("D1").Value = "ppm"

$ws.Columns.Item(9).Delete()

Write-Output "done"
